# Adds newly-scraped schedule rows for line 141 (update run at 30/12/2025 20:47:xx)
# across all three sheets, and refreshes the "last updated" / "total rows" header cells.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet "LP1912": append rows 632-647 --------------------------------
$ws1.Range("B632").Value = "20:47:31"
$ws1.Range("C632").Value = "20:52"
$ws1.Range("D632").Value = "15_ABASTO"
$ws1.Range("E632").Value = 5
$ws1.Range("F632").Value = "LP1912"
$ws1.Range("G632").Value = "30/12/2025"
$ws1.Range("B633").Value = "20:47:31"
$ws1.Range("C633").Value = "20:57"
$ws1.Range("D633").Value = "23_HERNANDEZ"
$ws1.Range("E633").Value = 10
$ws1.Range("F633").Value = "LP1912"
$ws1.Range("G633").Value = "30/12/2025"
$ws1.Range("B634").Value = "20:47:31"
$ws1.Range("C634").Value = "21:04"
$ws1.Range("D634").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("E634").Value = 17
$ws1.Range("F634").Value = "LP1912"
$ws1.Range("G634").Value = "30/12/2025"
$ws1.Range("B635").Value = "20:47:31"
$ws1.Range("C635").Value = "21:07"
$ws1.Range("D635").Value = "215B_EL PATO"
$ws1.Range("E635").Value = 20
$ws1.Range("F635").Value = "LP1912"
$ws1.Range("G635").Value = "30/12/2025"
$ws1.Range("B636").Value = "20:47:31"
$ws1.Range("C636").Value = "21:14"
$ws1.Range("D636").Value = "23_HERNANDEZ"
$ws1.Range("E636").Value = 27
$ws1.Range("F636").Value = "LP1912"
$ws1.Range("G636").Value = "30/12/2025"
$ws1.Range("B637").Value = "20:47:31"
$ws1.Range("C637").Value = "21:18"
$ws1.Range("D637").Value = "16_SANTA ANA"
$ws1.Range("E637").Value = 31
$ws1.Range("F637").Value = "LP1912"
$ws1.Range("G637").Value = "30/12/2025"
$ws1.Range("B638").Value = "20:47:31"
$ws1.Range("C638").Value = "21:21"
$ws1.Range("D638").Value = "26_HERNANDEZ"
$ws1.Range("E638").Value = 34
$ws1.Range("F638").Value = "LP1912"
$ws1.Range("G638").Value = "30/12/2025"
$ws1.Range("B639").Value = "20:47:31"
$ws1.Range("C639").Value = "21:23"
$ws1.Range("D639").Value = "15_ABASTO"
$ws1.Range("E639").Value = 36
$ws1.Range("F639").Value = "LP1912"
$ws1.Range("G639").Value = "30/12/2025"
$ws1.Range("B640").Value = "20:47:31"
$ws1.Range("C640").Value = "21:32"
$ws1.Range("D640").Value = "16_SANTA ANA"
$ws1.Range("E640").Value = 45
$ws1.Range("F640").Value = "LP1912"
$ws1.Range("G640").Value = "30/12/2025"
$ws1.Range("B641").Value = "20:47:31"
$ws1.Range("C641").Value = "21:32"
$ws1.Range("D641").Value = "23_HERNANDEZ"
$ws1.Range("E641").Value = 45
$ws1.Range("F641").Value = "LP1912"
$ws1.Range("G641").Value = "30/12/2025"
$ws1.Range("B642").Value = "20:47:31"
$ws1.Range("C642").Value = "21:38"
$ws1.Range("D642").Value = "17_ROMERO"
$ws1.Range("E642").Value = 51
$ws1.Range("F642").Value = "LP1912"
$ws1.Range("G642").Value = "30/12/2025"
$ws1.Range("B643").Value = "20:47:31"
$ws1.Range("C643").Value = "21:47"
$ws1.Range("D643").Value = "215A_EL PATO"
$ws1.Range("E643").Value = 60
$ws1.Range("F643").Value = "LP1912"
$ws1.Range("G643").Value = "30/12/2025"
$ws1.Range("B644").Value = "20:47:31"
$ws1.Range("C644").Value = "21:51"
$ws1.Range("D644").Value = "10_OLMOS"
$ws1.Range("E644").Value = 64
$ws1.Range("F644").Value = "LP1912"
$ws1.Range("G644").Value = "30/12/2025"
$ws1.Range("B645").Value = "20:47:31"
$ws1.Range("C645").Value = "22:08"
$ws1.Range("D645").Value = "17_ROMERO"
$ws1.Range("E645").Value = 81
$ws1.Range("F645").Value = "LP1912"
$ws1.Range("G645").Value = "30/12/2025"
$ws1.Range("B646").Value = "20:47:31"
$ws1.Range("C646").Value = "22:23"
$ws1.Range("D646").Value = "26_HERNANDEZ"
$ws1.Range("E646").Value = 96
$ws1.Range("F646").Value = "LP1912"
$ws1.Range("G646").Value = "30/12/2025"
$ws1.Range("B647").Value = "20:47:31"
$ws1.Range("C647").Value = "22:25"
$ws1.Range("D647").Value = "10_OLMOS"
$ws1.Range("E647").Value = 98
$ws1.Range("F647").Value = "LP1912"
$ws1.Range("G647").Value = "30/12/2025"

# --- Sheet "LP1912-215": append rows 48-49 ------------------------------
$ws2.Range("B48").Value = "30/12/2025"
$ws2.Range("C48").Value = "20:47:31"
$ws2.Range("D48").Value = "21:07"
$ws2.Range("E48").Value = "215B_EL PATO"
$ws2.Range("F48").Value = 20
$ws2.Range("G48").Value = "LP1912"
$ws2.Range("B49").Value = "30/12/2025"
$ws2.Range("C49").Value = "20:47:31"
$ws2.Range("D49").Value = "21:47"
$ws2.Range("E49").Value = "215A_EL PATO"
$ws2.Range("F49").Value = 60
$ws2.Range("G49").Value = "LP1912"

# --- Sheet "6203-6173": append rows 78-80 -------------------------------
$ws3.Range("B78").Value = "30/12/2025"
$ws3.Range("C78").Value = "20:47:37"
$ws3.Range("D78").Value = "21:29"
$ws3.Range("E78").Value = "215C_LA PLATA"
$ws3.Range("F78").Value = 42
$ws3.Range("G78").Value = "L6203"
$ws3.Range("B79").Value = "30/12/2025"
$ws3.Range("C79").Value = "20:47:42"
$ws3.Range("D79").Value = "22:04"
$ws3.Range("E79").Value = "215A_LA PLATA"
$ws3.Range("F79").Value = 77
$ws3.Range("G79").Value = "L6173"
$ws3.Range("B80").Value = "30/12/2025"
$ws3.Range("C80").Value = "20:47:42"
$ws3.Range("D80").Value = "22:20"
$ws3.Range("E80").Value = "215B_LP-P MOR-40 Y 115"
$ws3.Range("F80").Value = 93
$ws3.Range("G80").Value = "L6173"

# --- Refresh header cells (A2 "Ultima actualizacion", A3 "Total filas") --
$ws1.Range("A2").Value = "Última actualización: 30/12/2025 20:47:42"
$ws1.Range("A3").Value = "Total filas: 646"

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 20:47:42"
$ws2.Range("A3").Value = "Total filas: 48"

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 20:47:42"
$ws3.Range("A3").Value = "Total filas: 79"
